# Generate Report for Handback
# Marks the zh-cn and de-de localization rows as handed back: updates the
# Status column, stamps the Latest Target File / Latest Handback File /
# Latest Handback DateTime columns, and links the new Target File cell back
# to the source markdown file (mirroring column A's hyperlink).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$sourceFileName = "a7a83e1b-a1ef-4662-9dfa-71fbb3930458.md"
$sourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/cdedcf74dde14b90a2afd19531234f99b47ff9ab/e2e/a7a83e1b-a1ef-4662-9dfa-71fbb3930458.md"

# ---- Overview sheet: broaden the zh-cn / de-de status columns ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---- zh-cn handback ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Range("J2").Value = $sourceFileName
$zh.Hyperlinks.Add($zh.Range("J2"), $sourceUrl, "", "", $sourceFileName)
$zh.Range("K2").Value = "a7a83e1b-a1ef-4662-9dfa-71fbb3930458.9ce10655ab9552006c5972a28b1d2b9b6fa1d724.zh-cn.xlf"
$zh.Range("L2").Value = "2017-03-02 08:06:30"
$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(10).ColumnWidth = 40
$zh.Columns.Item(11).ColumnWidth = 40

# ---- de-de handback ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Range("J2").Value = $sourceFileName
$de.Hyperlinks.Add($de.Range("J2"), $sourceUrl, "", "", $sourceFileName)
$de.Range("K2").Value = "a7a83e1b-a1ef-4662-9dfa-71fbb3930458.9ce10655ab9552006c5972a28b1d2b9b6fa1d724.de-de.xlf"
$de.Range("L2").Value = "2017-03-02 08:06:53"
$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(10).ColumnWidth = 40
$de.Columns.Item(11).ColumnWidth = 40
